# Applies the scheduled market-data refresh: updates currentAveragePrice /
# LevePrice / LeveProfit columns (H:N) for the affected leve rows across all
# eight job sheets, per the upstream runner diff.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# row 18
$ws.Range("H18").Value = 4984.75
$ws.Range("I18").Value = 4984.75
$ws.Range("K18").Value = 4984.75
$ws.Range("M18").Value = -4700.75
# row 41
$ws.Range("H41").Value = 1666.875
$ws.Range("I41").Value = 1367
$ws.Range("J41").Value = 1766.8334
$ws.Range("K41").Value = 1367
$ws.Range("L41").Value = 1766.8334
$ws.Range("M41").Value = -927
$ws.Range("N41").Value = -2646.8334
# row 51
$ws.Range("H51").Value = 125375210
$ws.Range("I51").Value = 500273
$ws.Range("J51").Value = 500000000
$ws.Range("K51").Value = 500273
$ws.Range("L51").Value = 500000000
$ws.Range("M51").Value = -499789
$ws.Range("N51").Value = -500000968
# row 86
$ws.Range("H86").Value = 58828300
$ws.Range("I86").Value = 71433144
$ws.Range("K86").Value = 71433144
$ws.Range("M86").Value = -71432021
# row 89
$ws.Range("H89").Value = 58828300
$ws.Range("I89").Value = 71433144
$ws.Range("K89").Value = 357165720
$ws.Range("M89").Value = -357160104
# row 107
$ws.Range("H107").Value = 683.96155
$ws.Range("I107").Value = 706.6667
$ws.Range("J107").Value = 632.875
$ws.Range("K107").Value = 706.6667
$ws.Range("L107").Value = 632.875
$ws.Range("M107").Value = 1213.3333
$ws.Range("N107").Value = -4472.875
# row 135
$ws.Range("H135").Value = 992.64703
$ws.Range("I135").Value = 1092.8
$ws.Range("K135").Value = 9835.199999999999
$ws.Range("M135").Value = -7300.199999999999
# row 138
$ws.Range("H138").Value = 5268.651
$ws.Range("I138").Value = 7616
$ws.Range("J138").Value = 4255.023
$ws.Range("K138").Value = 22848
$ws.Range("L138").Value = 12765.069
$ws.Range("M138").Value = -17708
$ws.Range("N138").Value = -23045.069
# row 141
$ws.Range("H141").Value = 2926.5454
$ws.Range("I141").Value = 2531.889
$ws.Range("J141").Value = 4702.5
$ws.Range("K141").Value = 7595.667
$ws.Range("L141").Value = 14107.5
$ws.Range("M141").Value = -2415.667
$ws.Range("N141").Value = -24467.5
$ws = $wb.Worksheets.Item("ARM")
# row 32
$ws.Range("H32").Value = 134250.55
$ws.Range("I32").Value = 139337.19
$ws.Range("K32").Value = 139337.19
$ws.Range("M32").Value = -139050.19
# row 45
$ws.Range("H45").Value = 402497
$ws.Range("I45").Value = 253121.75
$ws.Range("K45").Value = 253121.75
$ws.Range("M45").Value = -252744.75
# row 74
$ws.Range("H74").Value = 1013240.5
$ws.Range("I74").Value = 1991.9375
$ws.Range("J74").Value = 3709903.2
$ws.Range("K74").Value = 1991.9375
$ws.Range("L74").Value = 3709903.2
$ws.Range("M74").Value = -1117.9375
$ws.Range("N74").Value = -3711651.2
# row 77
$ws.Range("H77").Value = 1013240.5
$ws.Range("I77").Value = 1991.9375
$ws.Range("J77").Value = 3709903.2
$ws.Range("K77").Value = 9959.6875
$ws.Range("L77").Value = 18549516
$ws.Range("M77").Value = -5591.6875
$ws.Range("N77").Value = -18558252
# row 102
$ws.Range("H102").Value = 2161.6
$ws.Range("I102").Value = 2063.5
$ws.Range("K102").Value = 2063.5
$ws.Range("M102").Value = -441.5
# row 132
$ws.Range("H132").Value = 3337.7144
$ws.Range("I132").Value = 0
$ws.Range("J132").Value = 3337.7144
$ws.Range("K132").Value = 0
$ws.Range("L132").Value = 10013.1432
$ws.Range("M132").ClearContents()
$ws.Range("N132").Value = -15073.1432
$ws = $wb.Worksheets.Item("BSM")
# row 20
$ws.Range("H20").Value = 903.6923
$ws.Range("I20").Value = 842.8823
$ws.Range("J20").Value = 1018.55554
$ws.Range("K20").Value = 842.8823
$ws.Range("L20").Value = 1018.55554
$ws.Range("M20").Value = -595.8823
$ws.Range("N20").Value = -1512.55554
# row 105
$ws.Range("H105").Value = 4768.1562
$ws.Range("I105").Value = 4921.76
$ws.Range("K105").Value = 4921.76
$ws.Range("M105").Value = -3174.76
# row 107
$ws.Range("H107").Value = 6680.8335
$ws.Range("I107").Value = 7675.184
$ws.Range("K107").Value = 7675.184
$ws.Range("M107").Value = -5755.184
# row 122
$ws.Range("H122").Value = 49999
$ws.Range("J122").Value = 49999
$ws.Range("L122").Value = 49999
$ws.Range("N122").Value = -59799
$ws = $wb.Worksheets.Item("CRP")
# row 16
$ws.Range("H16").Value = 3171
$ws.Range("I16").Value = 2474.9
$ws.Range("K16").Value = 2474.9
$ws.Range("M16").Value = -2187.9
# row 31
$ws.Range("H31").Value = 2904.1965
$ws.Range("I31").Value = 2834.6
$ws.Range("K31").Value = 2834.6
$ws.Range("M31").Value = -2539.6
# row 34
$ws.Range("H34").Value = 2904.1965
$ws.Range("I34").Value = 2834.6
$ws.Range("K34").Value = 2834.6
$ws.Range("M34").Value = -2632.6
# row 113
$ws.Range("H113").Value = 3171
$ws.Range("I113").Value = 2474.9
$ws.Range("K113").Value = 2474.9
$ws.Range("M113").Value = -304.9000000000001
# row 132
$ws.Range("H132").Value = 19935.053
$ws.Range("I132").Value = 25228.814
$ws.Range("K132").Value = 75686.442
$ws.Range("M132").Value = -73156.442
$ws = $wb.Worksheets.Item("CUL")
# row 14
$ws.Range("H14").Value = 95.14286
$ws.Range("I14").Value = 95.14286
$ws.Range("K14").Value = 285.42858
$ws.Range("M14").Value = -112.42858
# row 33
$ws.Range("H33").Value = 20100466
$ws.Range("J33").Value = 28714642
$ws.Range("L33").Value = 172287852
$ws.Range("N33").Value = -172288418
# row 86
$ws.Range("H86").Value = 1954.7778
$ws.Range("J86").Value = 1636.625
$ws.Range("L86").Value = 4909.875
$ws.Range("N86").Value = -7281.875
# row 89
$ws.Range("H89").Value = 1954.7778
$ws.Range("J89").Value = 1636.625
$ws.Range("L89").Value = 14729.625
$ws.Range("N89").Value = -26585.625
$ws = $wb.Worksheets.Item("GSM")
# row 80
$ws.Range("H80").Value = 100632.664
$ws.Range("J80").Value = 55499.715
$ws.Range("L80").Value = 55499.715
$ws.Range("N80").Value = -57495.715
# row 83
$ws.Range("H83").Value = 100632.664
$ws.Range("J83").Value = 55499.715
$ws.Range("L83").Value = 277498.575
$ws.Range("N83").Value = -287482.575
# row 93
$ws.Range("H93").Value = 100000
$ws.Range("I93").Value = 100000
$ws.Range("K93").Value = 100000
$ws.Range("M93").Value = -98128
# row 97
$ws.Range("H97").Value = 3352.1428
$ws.Range("J97").Value = 11686.667
$ws.Range("L97").Value = 11686.667
$ws.Range("N97").Value = -12678.667
# row 102
$ws.Range("H102").Value = 33334642
$ws.Range("I102").Value = 41668052
$ws.Range("J102").Value = 1000
$ws.Range("K102").Value = 41668052
$ws.Range("L102").Value = 1000
$ws.Range("M102").Value = -41666430
$ws.Range("N102").Value = -4244
# row 111
$ws.Range("H111").Value = 39997.5
$ws.Range("J111").Value = 39997.5
$ws.Range("L111").Value = 39997.5
$ws.Range("N111").Value = -46131.5
# row 122
$ws.Range("H122").Value = 1778.7084
$ws.Range("I122").Value = 1457.0286
$ws.Range("K122").Value = 4371.085800000001
$ws.Range("M122").Value = -1921.085800000001
# row 126
$ws.Range("H126").Value = 3168.9
$ws.Range("I126").Value = 2648.3333
$ws.Range("K126").Value = 7944.999899999999
$ws.Range("M126").Value = -5474.999899999999
$ws = $wb.Worksheets.Item("LTW")
# row 16
$ws.Range("H16").Value = 1587.7667
$ws.Range("I16").Value = 1415.5
$ws.Range("J16").Value = 3999.5
$ws.Range("K16").Value = 1415.5
$ws.Range("L16").Value = 3999.5
$ws.Range("M16").Value = -1245.5
$ws.Range("N16").Value = -4339.5
# row 40
$ws.Range("H40").Value = 2509.7104
$ws.Range("I40").Value = 1562.8148
$ws.Range("K40").Value = 1562.8148
$ws.Range("M40").Value = -1426.8148
# row 47
$ws.Range("H47").Value = 15000
$ws.Range("J47").Value = 15000
$ws.Range("L47").Value = 15000
$ws.Range("N47").Value = -15980
# row 52
$ws.Range("H52").Value = 15000
$ws.Range("J52").Value = 15000
$ws.Range("L52").Value = 15000
$ws.Range("N52").Value = -15466
# row 82
$ws.Range("H82").Value = 515.0833
$ws.Range("I82").Value = 539.1818
$ws.Range("K82").Value = 539.1818
$ws.Range("M82").Value = -178.1818
# row 85
$ws.Range("H85").Value = 515.0833
$ws.Range("I85").Value = 539.1818
$ws.Range("K85").Value = 539.1818
$ws.Range("M85").Value = 708.8182
# row 132
$ws.Range("H132").Value = 2477.077
$ws.Range("I132").Value = 2697.2144
$ws.Range("K132").Value = 8091.6432
$ws.Range("M132").Value = -5561.6432
# row 135
$ws.Range("H135").Value = 145000
$ws.Range("J135").Value = 145000
$ws.Range("L135").Value = 145000
$ws.Range("N135").Value = -155140
# row 136
$ws.Range("H136").Value = 3007.2766
$ws.Range("I136").Value = 2156.0715
$ws.Range("K136").Value = 6468.2145
$ws.Range("M136").Value = -3918.2145
$ws = $wb.Worksheets.Item("WVR")
# row 132
$ws.Range("H132").Value = 2127.639
$ws.Range("I132").Value = 1181.16
$ws.Range("K132").Value = 3543.48
$ws.Range("M132").Value = -1013.48
